$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: remove the stray "_GoBack" bookmark that currently sits
# right after the due-date year ("...2018") in the second paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# Change 2: in the "Stretch Levels" paragraph, replace "extra credit"
# with "a reputation bonus" and drop a fresh "_GoBack" bookmark right
# after the inserted phrase (mirroring how Word marks the last edit).
# Temporary bookmarks are used purely to force the run to split at the
# exact boundaries, then are removed again once the split has taken
# hold.
# ------------------------------------------------------------------
$boundary1 = $d.Content
$boundary1.Find.Execute(", try to complete these stretch levels for extra credit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitBeforeRun = $d.Range($boundary1.Start, $boundary1.Start)
$d.Bookmarks.Add("TempSplitBefore", $splitBeforeRun)

$creditRange = $d.Content
$creditRange.Find.Execute("extra credit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitBeforeBonus = $d.Range($creditRange.Start, $creditRange.Start)
$d.Bookmarks.Add("TempSplitMid", $splitBeforeBonus)
$splitAfterBonus = $d.Range($creditRange.End, $creditRange.End)
$d.Bookmarks.Add("TempSplitAfter", $splitAfterBonus)

$replaceRange = $d.Content
$replaceRange.Find.Execute("extra credit", $true, $false, $false, $false, $false, $true, 1, $false, "a reputation bonus", 2) | Out-Null

$goBackRange = $d.Bookmarks("TempSplitAfter").Range
$d.Bookmarks.Add("_GoBack", $goBackRange)

$d.Bookmarks("TempSplitBefore").Delete()
$d.Bookmarks("TempSplitMid").Delete()
$d.Bookmarks("TempSplitAfter").Delete()

# ------------------------------------------------------------------
# Change 3: reword rule 3 of the submission rules list.
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("All external ", $true, $false, $false, $false, $false, $true, 1, $false, "Any resources not created by you (images, ", 2) | Out-Null

$r2 = $d.Content
$r2.Find.Execute(" libraries must be referenced using a CDN, not directly included in your assignment submission.", $true, $false, $false, $false, $false, $true, 1, $false, " libraries, etc.) must be referenced using a CDN or URL, not directly included in your assignment submission.", 2) | Out-Null
